# Apply the edits described by the commit:
# "Support to convert domain name URL to IP address"
#
# Summary of changes to the first worksheet ("咕咕咕"):
#  - B7  : http://ah.ssi.996icu.net:8082        -> http://www.bilibili.com:8082
#  - B8  : ah.hahah.com                          -> www.acfun.cn  (+ new hyperlink)
#  - B9  : http://192.168.50.234:8088/os/abc     -> two lines (adds a second IP:port URL),
#          wrapped text, taller row, hyperlink display text kept as first URL
#  - B10 : https://123.256.23.5:4433/23index.html -> https://123.251.23.5:4433/23index.html
#  - selection / view moves to C10 with the sheet scrolled down a bit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# NOTE: indexing into $ws.Hyperlinks with .Item(n) returns a proxy whose
# properties (.Address/.TextToDisplay/.Delete) don't bind correctly in this
# runtime. Collecting references via a foreach loop first works reliably,
# so we grab all the existing hyperlinks up front (in B2,B3,B4,B6,B7,B9,B10,B12
# order) and mutate those captured references instead.
$existingLinks = @()
foreach ($h in $ws.Hyperlinks) { $existingLinks += $h }

$hlB7  = $existingLinks[4]   # B7
$hlB9  = $existingLinks[5]   # B9
$hlB10 = $existingLinks[6]   # B10

# --- B7: update the URL text and its hyperlink target -------------------
$ws.Range("B7").Value = "http://www.bilibili.com:8082"
$hlB7.Address = "http://www.bilibili.com:8082"

# --- B8: update the domain text and add a brand new hyperlink -----------
$ws.Range("B8").Value = "www.acfun.cn"
$ws.Hyperlinks.Add($ws.Range("B8"), "http://www.acfun.cn") | Out-Null

# --- B9: add a second URL on a new line, wrap text, taller row ----------
$firstUrl = "http://192.168.50.234:8088/os/abc"
$secondUrl = "http://182.78.76.234:8081/os/abc"
$ws.Range("B9").Value = $firstUrl + $nl + $secondUrl
$ws.Range("B9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 27.75
$hlB9.Address = $firstUrl
$hlB9.TextToDisplay = $firstUrl

# --- B10: fix the IP typo in the URL -------------------------------------
$ws.Range("B10").Value = "https://123.251.23.5:4433/23index.html"
$hlB10.Address = "https://123.251.23.5:4433/23index.html"

# --- Update the view: scroll down and select C10 -------------------------
$ws.Activate()
$ws.Range("C10").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 8
